$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "00:00:00:114"
$ws.Range("D4").Value = "00:00:00:013"
$ws.Range("D5").Value = "00:00:00:013"
$ws.Range("D6").Value = "00:00:00:013"
$ws.Range("D7").Value = "00:00:00:013"
$ws.Range("D8").Value = "00:00:00:013"
$ws.Range("D9").Value = "00:00:00:013"
$ws.Range("D10").Value = "00:00:00:013"
$ws.Range("D11").Value = "00:00:00:013"
$ws.Range("D12").Value = "00:00:00:012"
$ws.Range("D13").Value = "00:00:00:013"
$ws.Range("D14").Value = "00:00:00:013"
$ws.Range("D15").Value = "00:00:00:013"
$ws.Range("D16").Value = "00:00:00:013"
$ws.Range("D17").Value = "00:00:00:013"
$ws.Range("D18").Value = "00:00:00:013"
$ws.Range("D19").Value = "00:00:00:013"
$ws.Range("D20").Value = "00:00:00:013"
$ws.Range("D21").Value = "00:00:00:013"
$ws.Range("D22").Value = "00:00:00:012"
$ws.Range("D23").Value = "00:00:00:013"
$ws.Range("D24").Value = "00:00:00:013"
$ws.Range("D25").Value = "00:00:00:013"
$ws.Range("D26").Value = "00:00:00:013"
$ws.Range("D27").Value = "00:00:00:013"
$ws.Range("D28").Value = "00:00:00:013"
$ws.Range("D29").Value = "00:00:00:013"
$ws.Range("D30").Value = "00:00:00:013"
$ws.Range("D31").Value = "00:00:00:012"
$ws.Range("D32").Value = "00:00:00:012"
$ws.Range("D33").Value = "00:00:00:012"
$ws.Range("D34").Value = "00:00:00:013"
$ws.Range("D35").Value = "00:00:00:016"
$ws.Range("D36").Value = "00:00:00:013"
$ws.Range("D37").Value = "00:00:00:012"
$ws.Range("D38").Value = "00:00:00:013"
$ws.Range("D39").Value = "00:00:00:012"
$ws.Range("D40").Value = "00:00:00:013"
$ws.Range("D41").Value = "00:00:00:013"
$ws.Range("D42").Value = "00:00:00:012"
$ws.Range("D43").Value = "00:00:00:012"
$ws.Range("D44").Value = "00:00:00:012"
$ws.Range("D45").Value = "00:00:00:012"
$ws.Range("D46").Value = "00:00:00:013"
$ws.Range("D47").Value = "00:00:00:012"
$ws.Range("D48").Value = "00:00:00:013"
$ws.Range("D49").Value = "00:00:00:012"
$ws.Range("D50").Value = "00:00:00:012"
$ws.Range("D51").Value = "00:00:00:013"
$ws.Range("D52").Value = "00:00:00:013"
$ws.Range("D53").Value = "00:00:00:017"
$ws.Range("D54").Value = "00:00:00:013"
$ws.Range("D55").Value = "00:00:00:013"
$ws.Range("D56").Value = "00:00:00:013"
$ws.Range("D57").Value = "00:00:00:013"
$ws.Range("D58").Value = "00:00:00:013"
$ws.Range("D59").Value = "00:00:00:013"
$ws.Range("D60").Value = "00:00:00:025"
$ws.Range("D61").Value = "00:00:00:013"
$ws.Range("D62").Value = "00:00:00:013"
$ws.Range("D63").Value = "00:00:00:013"
$ws.Range("D64").Value = "00:00:00:013"
